$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price/Volume columns so numeric-looking strings stay as text
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "42.063.76"
$ws.Range("E2").Value = "  -0.78%  "

$ws.Range("D3").Value = "2.249.05"
$ws.Range("E3").Value = "  -1.05%  "

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "306.79"
$ws.Range("E5").Value = "  +0.02%  "

$ws.Range("D6").Value = "96.72"
$ws.Range("E6").Value = "  -0.92%  "

$ws.Range("D7").Value = "0.523"
$ws.Range("E7").Value = "  -1.31%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").Value = "0.489"
$ws.Range("E9").Value = "  -0.50%  "

$ws.Range("D10").Value = "34.94"
$ws.Range("E10").Value = "  -2.03%  "

$ws.Range("D11").Value = "0.0816"
$ws.Range("E11").Value = "  +2.58%  "

$ws.Range("E12").Value = "  +1.27%  "

$ws.Range("D13").Value = "6.78"
$ws.Range("E13").Value = "  +1.69%  "

$ws.Range("D14").Value = "2.594.92"
$ws.Range("E14").Value = "  -1.14%  "

$ws.Range("D15").Value = "14.44"
$ws.Range("E15").Value = "  +0.49%  "

$ws.Range("D16").Value = "2.242.91"
$ws.Range("E16").Value = "  -0.33%  "

$ws.Range("D17").Value = "0.781"
$ws.Range("E17").Value = "  -1.70%  "

$ws.Range("D18").Value = "41.924.42"
$ws.Range("E18").Value = "  -0.83%  "

$ws.Range("D19").Value = "12.20"
$ws.Range("E19").Value = "  -2.46%  "

$ws.Range("D20").Value = "0.0₃0902"
$ws.Range("E20").Value = "  -0.85%  "

$ws.Range("D21").Value = "5.92"
$ws.Range("E21").Value = "  -0.49%  "

$ws.Range("D22").Value = "67.11"
$ws.Range("E22").Value = "  -0.70%  "

$ws.Range("D23").Value = "235.78"
$ws.Range("E23").Value = "  -1.95%  "

$ws.Range("D24").Value = "2.58"
$ws.Range("E24").Value = "  -0.52%  "

$ws.Range("D25").Value = "1.94"
$ws.Range("E25").Value = "  -0.59%  "

$ws.Range("D26").Value = "1.01"
$ws.Range("E26").Value = "  +0.92%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "23.30"
$ws.Range("E27").Value = "  -2.23%  "

$ws.Range("B28").Value = "InjectiveProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D28").Value = "37.76"
$ws.Range("E28").Value = "  +0.82%  "

$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").Value = "9.47"
$ws.Range("E29").Value = "  -0.33%  "

$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "2.11"
$ws.Range("E30").Value = "  +0.03%  "

$ws.Range("D31").Value = "167.84"
$ws.Range("E31").Value = "  +4.94%  "

$ws.Range("E32").Value = "  -0.02%  "

$ws.Range("D33").Value = "5.17"
$ws.Range("E33").Value = "  -1.74%  "

$ws.Range("B34").Value = "Celestia"
$ws.Range("C34").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D34").Value = "17.52"
$ws.Range("E34").Value = "  +2.62%  "

$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").Value = "3.04"
$ws.Range("E35").Value = "  -3.46%  "

$ws.Range("D36").Value = "0.0719"
$ws.Range("E36").Value = "  -2.98%  "

$ws.Range("E37").Value = "  +1.49%  "

$ws.Range("E38").Value = "  -0.01%  "

$ws.Range("E39").Value = "  -2.81%  "

$ws.Range("D40").Value = "1.79"
$ws.Range("E40").Value = "  -2.37%  "

$ws.Range("D41").Value = "4.06"
$ws.Range("E41").Value = "  -0.28%  "

$ws.Range("D42").Value = "1.939.61"
$ws.Range("E42").Value = "  -2.93%  "

$ws.Range("D43").Value = "0.0281"
$ws.Range("E43").Value = "  -1.38%  "

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "18.49"
$ws.Range("E44").Value = "  -2.23%  "

$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "2.17"
$ws.Range("E45").Value = "  -11.09%  "

$ws.Range("D46").Value = "2.89"
$ws.Range("E46").Value = "  -1.72%  "

$ws.Range("D47").Value = "9.64"
$ws.Range("E47").Value = "  -3.41%  "

$ws.Range("D48").Value = "53.99"
$ws.Range("E48").Value = "  +1.75%  "

$ws.Range("D49").Value = "2.465.42"
$ws.Range("E49").Value = "  -1.16%  "

$ws.Range("D50").Value = "71.13"
$ws.Range("E50").Value = "  -1.49%  "

$ws.Range("D51").Value = "91.09"
$ws.Range("E51").Value = "  -0.50%  "
